# Adds two new paragraphs right after the "No conformidad en la pegada"
# paragraph: an empty paragraph, followed by a paragraph containing the
# text "15:45 domingo agregamos algo mas". The trailing "_GoBack" bookmark
# (originally sitting right after "pegada") ends up positioned right after
# the newly typed text, exactly as in the target XML.

$d = $word.ActiveDocument

$anchorText  = "No conformidad en la pegada"
$newLine     = "15:45 domingo agregamos algo mas"
$tail1       = "ZzPlaceholderEmptyParaZz"
$tail2       = "ZzPlaceholderBookmarkZz"

# 1) The original "_GoBack" bookmark sits right at the end of $anchorText.
#    Drop it now; it gets re-created in its new location once the new
#    text exists (empty/collapsed Bookmarks.Add at a "paragraph end"
#    position behaves oddly in this host, so we add it while it is still
#    mid-paragraph -- see step 3).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Locate the end of the anchor paragraph's text and type the new
#    content there: paragraph break, a throwaway marker (so the freshly
#    created empty paragraph actually contains a real character while we
#    work -- avoids the host leaving a stray empty <w:r/> behind once we
#    delete the marker), another paragraph break, the real new sentence,
#    and a second throwaway marker right after it (keeps the insertion
#    point away from the "end of paragraph" position until we are done).
$rng = $d.Content
$rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.Select()
$word.Selection.TypeParagraph() | Out-Null
$word.Selection.TypeText($tail1) | Out-Null
$word.Selection.TypeParagraph() | Out-Null
$word.Selection.TypeText($newLine + $tail2) | Out-Null

# 3) Re-add "_GoBack" right after the real text (before $tail2), while
#    that position is still mid-paragraph.
$rngMark = $d.Content
$rngMark.Find.Execute($newLine, $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0) | Out-Null
$rngMark.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rngMark)

# 4) Remove the two throwaway markers now that everything is positioned.
$rngTail2 = $d.Content
$rngTail2.Find.Execute($tail2, $true, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0) | Out-Null
$rngTail2.Delete()

$rngTail1 = $d.Content
$rngTail1.Find.Execute($tail1, $true, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0) | Out-Null
$rngTail1.Delete()
